$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.009.50"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "1.774.52"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "328.54"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "0.4505"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "0.3562"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").Value = "0.07445"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "41.96"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("D11").Value = "1.109"
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("D12").Value = "0.9990"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "21.02"
$ws.Range("D14").Value = "6.055"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").Value = "7.263"
$ws.Range("E15").Value = "  +2.82%  "
$ws.Range("D16").Value = "1.773.05"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").Value = "94.05"
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("D18").Value = "0.00001066"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").Value = "0.06448"
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("D20").Value = "0.9990"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("D22").Value = "5.795"
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("D23").Value = "28.007.37"
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("D24").Value = "11.33"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("D25").Value = "2.129"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").Value = "162.02"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "20.42"
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("D28").Value = "1.975.33"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("D29").Value = "2.169"
$ws.Range("E29").Value = "  +7.02%  "
$ws.Range("D30").Value = "124.75"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "1.113"
$ws.Range("E31").Value = "  +6.43%  "
$ws.Range("D32").Value = "5.721"
$ws.Range("E32").Value = "  +6.50%  "
$ws.Range("D33").Value = "0.09226"
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("D34").Value = "3.690"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").Value = "11.89"
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("D36").Value = "0.06218"
$ws.Range("E36").Value = "  +3.78%  "
$ws.Range("D37").Value = "0.02294"
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").Value = "0.2115"
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("D39").Value = "4.996"
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("D40").Value = "0.6329"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").Value = "1.189"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").Value = "7.917"
$ws.Range("E43").Value = "  +2.95%  "
$ws.Range("D44").Value = "13.32"
$ws.Range("E44").Value = "  +1.78%  "
$ws.Range("D45").Value = "3.757"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("D46").Value = "0.5903"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("D47").Value = "122.87"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "1.964"
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("D49").Value = "1.143"
$ws.Range("D50").Value = "0.06895"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").Value = "73.03"
$ws.Range("E51").Value = "  +2.60%  "
